# The BF column holds a "Date" string that was off by one day due to how the
# NBA stats site displayed dates (season label "5-20-2007-08" really means
# the game date 2008-05-20). Fix it for every data row (BF2:BF31).
#
# NumberFormat is temporarily set to Text ("@") before the assignment so the
# date-looking literal "2008-05-20" is stored as a plain string instead of
# being auto-converted to a date serial number, then ClearFormats() restores
# the cells' original (default/general) formatting so no visible formatting
# changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "2008-05-20"
}

$dateRange.ClearFormats()
